$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Decrement the "id" column (A2:A5) by 1, making it 0-indexed
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3

# Highlight the current section on the timeline by moving the active cell selection
$ws.Range("C11").Select()
